$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F (copy style from existing header cell B1)
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Timestamps (as strings, matching the inlineStr type in the diff) for rows 2-14
$timestamps = @(
    "2021-10-05 10:51:39.357129",
    "2021-10-05 10:51:39.357139",
    "2021-10-05 10:51:39.357142",
    "2021-10-05 10:51:39.357145",
    "2021-10-05 10:51:39.357148",
    "2021-10-05 10:51:39.357150",
    "2021-10-05 10:51:39.357153",
    "2021-10-05 10:51:39.357155",
    "2021-10-05 10:51:39.357158",
    "2021-10-05 10:51:39.357161",
    "2021-10-05 10:51:39.357163",
    "2021-10-05 10:51:39.357166",
    "2021-10-05 10:51:39.357168"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
